# Fill in the missing "CODIGO" (column F) values for rows 13-111 of the
# "caña" sheet with CAN-012 .. CAN-110, copy the existing cell formatting
# from F12 down, turn off the AutoFilter, and update the sheet's selection
# / scroll position to match the saved view left by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 13; $row -le 111; $row++) {
    $num = $row - 1
    $code = "CAN-{0:D3}" -f $num
    $ws.Cells.Item($row, 6).Value = $code
}

# Match the formatting (style, border) already used by F2:F12 in the same column.
$ws.Range("F12").Copy()
$ws.Range("F13:F111").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The table no longer keeps its AutoFilter dropdowns.
$ws.AutoFilterMode = $false

# Leave the sheet scrolled/selected the way the author saved it.
$ws.Activate()
$ws.Range("F12:F111").Select()
$excel.ActiveWindow.ScrollRow = 85
$excel.ActiveWindow.ScrollColumn = 5
